$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '23.711.28'
$ws.Range('E2').Value = '  +1.26%  '
Set-TextValue $ws.Range('D3') '1.652.83'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('E5').Value = '  +0.10%  '
Set-TextValue $ws.Range('D6') '303.20'
$ws.Range('E6').Value = '  -0.08%  '
Set-TextValue $ws.Range('D7') '0.3801'
$ws.Range('E7').Value = '  +0.57%  '
Set-TextValue $ws.Range('D8') '0.3619'
$ws.Range('E8').Value = '  -0.07%  '
Set-TextValue $ws.Range('D9') '50.99'
$ws.Range('E9').Value = '  -1.96%  '
Set-TextValue $ws.Range('D10') '1.247'
$ws.Range('E10').Value = '  +1.82%  '
Set-TextValue $ws.Range('D11') '0.08203'
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('E12').Value = '  +0.14%  '
Set-TextValue $ws.Range('D13') '22.67'
$ws.Range('E13').Value = '  +1.49%  '
Set-TextValue $ws.Range('D14') '6.520'
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('E15').Value = '  +1.16%  '
Set-TextValue $ws.Range('D16') '0.00001234'
$ws.Range('E16').Value = '  -0.47%  '
Set-TextValue $ws.Range('D17') '1.656.13'
$ws.Range('E17').Value = '  +1.73%  '
Set-TextValue $ws.Range('D18') '97.27'
$ws.Range('E18').Value = '  +2.53%  '
Set-TextValue $ws.Range('D19') '0.07018'
$ws.Range('E19').Value = '  +1.16%  '
Set-TextValue $ws.Range('D20') '6.785'
$ws.Range('E20').Value = '  +3.32%  '
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('E22').Value = '  +0.20%  '
Set-TextValue $ws.Range('D23') '12.86'
$ws.Range('E23').Value = '  +2.58%  '
Set-TextValue $ws.Range('D24') '23.721.93'
$ws.Range('E24').Value = '  +1.27%  '
Set-TextValue $ws.Range('D25') '2.530'
$ws.Range('E25').Value = '  +1.09%  '
Set-TextValue $ws.Range('D26') '3.052'
$ws.Range('E26').Value = '  -0.12%  '
Set-TextValue $ws.Range('D27') '21.27'
$ws.Range('E27').Value = '  +0.77%  '
Set-TextValue $ws.Range('D28') '151.50'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').Value = '  -1.12%  '
Set-TextValue $ws.Range('D30') '134.38'
$ws.Range('E30').Value = '  +0.98%  '
Set-TextValue $ws.Range('D31') '1.837.29'
$ws.Range('E31').Value = '  +1.36%  '
Set-TextValue $ws.Range('D32') '6.925'
$ws.Range('E32').Value = '  +4.56%  '
Set-TextValue $ws.Range('D33') '2.226'
$ws.Range('E33').Value = '  +3.10%  '
Set-TextValue $ws.Range('D34') '1.074'
$ws.Range('E34').Value = '  +2.59%  '
Set-TextValue $ws.Range('D35') '11.72'
$ws.Range('E35').Value = '  +3.89%  '
Set-TextValue $ws.Range('D36') '0.02808'
$ws.Range('E36').Value = '  +1.91%  '
Set-TextValue $ws.Range('D37') '0.2511'
$ws.Range('E37').Value = '  +0.81%  '
Set-TextValue $ws.Range('D38') '0.08820'
$ws.Range('E38').Value = '  +0.50%  '
Set-TextValue $ws.Range('D39') '6.102'
$ws.Range('E39').Value = '  +1.74%  '
Set-TextValue $ws.Range('D40') '0.07078'
$ws.Range('E40').Value = '  -0.34%  '
Set-TextValue $ws.Range('D41') '13.01'
Set-TextValue $ws.Range('D42') '0.7028'
$ws.Range('E42').Value = '  +0.57%  '
Set-TextValue $ws.Range('D43') '1.335'
$ws.Range('E43').Value = '  -0.34%  '
Set-TextValue $ws.Range('D44') '15.97'
$ws.Range('E44').Value = '  +0.49%  '
Set-TextValue $ws.Range('D45') '0.6503'
$ws.Range('E45').Value = '  +0.30%  '
Set-TextValue $ws.Range('D46') '2.320'
$ws.Range('E46').Value = '  +2.21%  '
$ws.Range('E47').Value = '  +0.01%  '
Set-TextValue $ws.Range('D48') '3.960'
$ws.Range('E48').Value = '  -0.23%  '
Set-TextValue $ws.Range('D49') '0.07951'
$ws.Range('E49').Value = '  -0.26%  '
Set-TextValue $ws.Range('D50') '128.04'
$ws.Range('E50').Value = '  +1.14%  '
Set-TextValue $ws.Range('D51') '1.187'
$ws.Range('E51').Value = '  +0.13%  '
